$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the genre value in C2 from "Action" to "Drama"
$ws.Range("C2").Value = "Drama"

# Move the active selection from C2 to C3
$ws.Range("C3").Select()
